$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($cell, $value) {
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = "Normal"
}

Set-TextCell $ws.Range("D2") "41.834.46"
Set-TextCell $ws.Range("E2") "  +5.36%  "
Set-TextCell $ws.Range("D3") "2.230.61"
Set-TextCell $ws.Range("E3") "  +2.64%  "
Set-TextCell $ws.Range("E4") "  -0.03%  "
Set-TextCell $ws.Range("D5") "231.89"
Set-TextCell $ws.Range("E5") "  +2.26%  "
Set-TextCell $ws.Range("E6") "  +0.19%  "
Set-TextCell $ws.Range("D7") "61.59"
Set-TextCell $ws.Range("E7") "  -2.34%  "
Set-TextCell $ws.Range("E8") "  -0.01%  "
Set-TextCell $ws.Range("E9") "  +2.79%  "
Set-TextCell $ws.Range("D10") "59.20"
Set-TextCell $ws.Range("E10") "  +1.21%  "
Set-TextCell $ws.Range("D11") "0.0898"
Set-TextCell $ws.Range("E11") "  +5.58%  "
Set-TextCell $ws.Range("E12") "  -0.15%  "
Set-TextCell $ws.Range("D13") "2.558.96"
Set-TextCell $ws.Range("E13") "  +2.60%  "
Set-TextCell $ws.Range("E14") "  -1.58%  "
Set-TextCell $ws.Range("D15") "22.04"
Set-TextCell $ws.Range("E15") "  +1.07%  "
Set-TextCell $ws.Range("D16") "0.803"
Set-TextCell $ws.Range("E16") "  -1.06%  "
Set-TextCell $ws.Range("D17") "5.59"
Set-TextCell $ws.Range("E17") "  +1.79%  "
Set-TextCell $ws.Range("D18") "2.219.12"
Set-TextCell $ws.Range("E18") "  +2.06%  "
Set-TextCell $ws.Range("D19") "41.693.17"
Set-TextCell $ws.Range("E19") "  +5.13%  "
Set-TextCell $ws.Range("D20") "0.0₃0901"
Set-TextCell $ws.Range("E20") "  -1.37%  "
Set-TextCell $ws.Range("D21") "72.08"
Set-TextCell $ws.Range("E21") "  +0.46%  "
Set-TextCell $ws.Range("D22") "6.03"
Set-TextCell $ws.Range("E22") "  +0.44%  "
Set-TextCell $ws.Range("D23") "250.18"
Set-TextCell $ws.Range("E23") "  +9.05%  "
Set-TextCell $ws.Range("B25") "PancakeSwap"
Set-TextCell $ws.Range("C25") "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
Set-TextCell $ws.Range("D25") "2.40"
Set-TextCell $ws.Range("E25") "  +1.97%  "
Set-TextCell $ws.Range("B26") "Toncoin"
Set-TextCell $ws.Range("C26") "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
Set-TextCell $ws.Range("D26") "2.32"
Set-TextCell $ws.Range("E26") "  -0.24%  "
Set-TextCell $ws.Range("D27") "9.61"
Set-TextCell $ws.Range("E27") "  +0.45%  "
Set-TextCell $ws.Range("B28") "Monero"
Set-TextCell $ws.Range("C28") "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
Set-TextCell $ws.Range("D28") "167.51"
Set-TextCell $ws.Range("E28") "  -2.01%  "
Set-TextCell $ws.Range("B29") "Kaspa"
Set-TextCell $ws.Range("C29") "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
Set-TextCell $ws.Range("D29") "0.141"
Set-TextCell $ws.Range("E29") "  +1.32%  "
Set-TextCell $ws.Range("D30") "20.02"
Set-TextCell $ws.Range("E30") "  +0.93%  "
Set-TextCell $ws.Range("E31") "  -3.10%  "
Set-TextCell $ws.Range("D32") "2.68"
Set-TextCell $ws.Range("E32") "  +0.15%  "
Set-TextCell $ws.Range("E33") "  -0.25%  "
Set-TextCell $ws.Range("D34") "5.02"
Set-TextCell $ws.Range("E34") "  +6.84%  "
Set-TextCell $ws.Range("D35") "4.68"
Set-TextCell $ws.Range("E35") "  +3.27%  "
Set-TextCell $ws.Range("D36") "0.0638"
Set-TextCell $ws.Range("E36") "  +3.42%  "
Set-TextCell $ws.Range("E37") "  -4.53%  "
Set-TextCell $ws.Range("D38") "3.65"
Set-TextCell $ws.Range("E38") "  -5.78%  "
Set-TextCell $ws.Range("E39") "  -1.23%  "
Set-TextCell $ws.Range("D40") "0.000253"
Set-TextCell $ws.Range("E40") "  +28.14%  "
Set-TextCell $ws.Range("E41") "  +0.05%  "
Set-TextCell $ws.Range("D42") "4.92"
Set-TextCell $ws.Range("E42") "  -0.05%  "
Set-TextCell $ws.Range("E43") "  +4.35%  "
Set-TextCell $ws.Range("D44") "8.58"
Set-TextCell $ws.Range("E44") "  +8.72%  "
Set-TextCell $ws.Range("B45") "Cronos"
Set-TextCell $ws.Range("C45") "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
Set-TextCell $ws.Range("D45") "0.0979"
Set-TextCell $ws.Range("E45") "  +6.30%  "
Set-TextCell $ws.Range("B46") "TrustWalletToken"
Set-TextCell $ws.Range("C46") "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
Set-TextCell $ws.Range("D46") "1.23"
Set-TextCell $ws.Range("E46") "  +1.05%  "
Set-TextCell $ws.Range("D47") "99.00"
Set-TextCell $ws.Range("E47") "  -3.52%  "
Set-TextCell $ws.Range("D48") "1.482.82"
Set-TextCell $ws.Range("E48") "  -1.92%  "
Set-TextCell $ws.Range("B49") "InjectiveProtocol"
Set-TextCell $ws.Range("C49") "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
Set-TextCell $ws.Range("D49") "16.55"
Set-TextCell $ws.Range("E49") "  -6.49%  "
Set-TextCell $ws.Range("B50") "HuobiToken"
Set-TextCell $ws.Range("C50") "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
Set-TextCell $ws.Range("D50") "2.82"
Set-TextCell $ws.Range("E50") "  +0.51%  "
Set-TextCell $ws.Range("D51") "52.51"
Set-TextCell $ws.Range("E51") "  +5.90%  "
